$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep a text/string type even when the value looks
    # like a number (e.g. "100.00"), without leaving a residual number
    # format style applied to the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 3: employee name anonymized (vacation entry for a deleted user)
$ws.Range("B3").Value = "test"
Set-TextValue $ws.Range("H3") ""
Set-TextValue $ws.Range("K3") "100.00"
Set-TextValue $ws.Range("L3") "0.00"
Set-TextValue $ws.Range("M3") "0.0"

# Row 4: employee name anonymized (vacation entry for a deleted user)
$ws.Range("B4").Value = "test"
Set-TextValue $ws.Range("H4") ""
Set-TextValue $ws.Range("K4") "100.00"
Set-TextValue $ws.Range("L4") "0.00"
Set-TextValue $ws.Range("M4") "0.0"

# Row 6: totals updated to match the new fuel/phone/fee figures
$ws.Range("K6").Value = 200
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("O6").Value = 220
